# The last week (rows 29-30, day 27) was previously cut short after two
# entries; add back the missing third working interval for that day
# (16:00-18:00) as row 31, pushing the blank separator / summary rows
# (previously 31-34) down by one to 32-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 31 - this shifts the old row 31 (blank separator)
# and the summary rows below it down by one, and copies row 30's
# formatting (styles/number formats) into the new row.
$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = 2014
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = 27
$ws.Range("D31").Value = 0.66666666666666663
$ws.Range("E31").Value = 0.75
$ws.Range("F31").Formula = "=(E31-D31)*24*60"
$ws.Range("G31").Formula = "=F31/60"

# Make sure the new data row matches the look of the rows above it.
$ws.Range("F31").Style = $ws.Range("F30").Style
$ws.Range("G31").Style = $ws.Range("G30").Style
$ws.Range("F31").NumberFormat = $ws.Range("F30").NumberFormat
$ws.Range("G31").NumberFormat = $ws.Range("G30").NumberFormat

$ws.Range("F31").Select()
